# Actualizacion automatica 2025-11-10 16:30:09
#
# "CUMPLIMIENTO MENSUAL" sheet refresh:
#   - Row 2 (OTROS) becomes GRIFERIAS with new figures.
#   - A new OTROS row is inserted at row 3 (old PORCELANATO figures there
#     are replaced).
#   - PORCELANATO moves down to row 4 with refreshed figures.
#   - The TOTAL row shifts down to row 5 with recomputed sums.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Row 2: GRIFERIAS (was OTROS) -- number formats already in place ---
$ws.Range("B2").Value = "GRIFERIAS"
$ws.Range("C2").Value = 86.41
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 86.41
$ws.Range("F2").Value = 0

# --- Row 3: OTROS (new content; cell formats already match) ---
$ws.Range("A3").Value = "OFICINA-CATAECSA"
$ws.Range("B3").Value = "OTROS"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 3241.06
$ws.Range("E3").Value = -3241.06
$ws.Range("F3").Value = 0

# --- Row 4: PORCELANATO (new row; A4 is brand new, B4 loses the old
#     bold/right-aligned TOTAL formatting, C4:F4 keep their existing
#     currency/percent formats) ---
$ws.Range("A4").Value = "OFICINA-CATAECSA"
$ws.Range("B4").ClearFormats()
$ws.Range("B4").Value = "PORCELANATO"
$ws.Range("C4").Value = 26000
$ws.Range("D4").Value = 8238.030000000001
$ws.Range("E4").Value = 17761.97
$ws.Range("F4").Value = 0.3168473076923077

# --- Row 5: TOTAL (brand new row, shifted down from row 4) ---
$ws.Range("B5").Value = "TOTAL"
$ws.Range("B5").HorizontalAlignment = -4152  # xlRight
$ws.Range("C5").Value = 26086.41
$ws.Range("C5").NumberFormat = '"$"#,##0.00'
$ws.Range("D5").Value = 11479.09
$ws.Range("D5").NumberFormat = '"$"#,##0.00'
$ws.Range("E5").Value = 14607.32
$ws.Range("E5").NumberFormat = '"$"#,##0.00'
$ws.Range("F5").Value = 0.4400410021923293
$ws.Range("F5").NumberFormat = "0.00%"

# --- Column width changes (D: 14 -> 13, F: 18 -> 24) ---
$ws.Columns.Item(4).ColumnWidth = 13 - 5/6
$ws.Columns.Item(6).ColumnWidth = 24 - 5/6
